$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "ejuice,e juice,vape liquid,vape juice,e liquid,eliquid"
$ws.Range("B6").Value = "vape pen"
$ws.Range("B7").Value = "vape pod, pod pystem, pod, pod mod "
$ws.Range("B8").Value = "disposable"
$ws.Range("B9").Value = "vape subscription, vape box "
$ws.Range("B10").Value = "hookah, cannabis, weed, thc, cbd, marijuana"
$ws.Range("B11").Value = "vape, vaping, smoke, vapor"

$ws.Range("C16").Select()
